# Auto-generated from the canonical OOXML diff (data/recommandations.xlsx).
# Rewrites the two sorted ranking tables cell-by-cell to their post-refresh
# values (BRVM automated daily market-data update), then drops the row that
# fell out of the Top-43 ranking on sheet "Recommandations".

$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Range("A1").Value = "Titre"
$ws1.Range("B1").Value = "Jours en Hausse"
$ws1.Range("C1").Value = "Jours en Baisse"
$ws1.Range("D1").Value = "Variation Totale (%)"
$ws1.Range("E1").Value = "Dernière Variation (%)"
$ws1.Range("F1").Value = "Recommandation"
$ws1.Range("G1").Value = "Stratégie"

$ws1.Range("A2").Value = "NEI-CEDA CI"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 3795
$ws1.Range("E2").Value = 945
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 8
$ws1.Range("D3").Value = 3375.88
$ws1.Range("E3").Value = 115.82
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "AIR LIQUIDE CI"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2765
$ws1.Range("E4").Value = 685
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2364.8
$ws1.Range("E5").Value = 580
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2004.99
$ws1.Range("E6").Value = 503.61
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "BRVM - TRANSPORT"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 1458.85
$ws1.Range("E7").Value = 364.71
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "BRVM - AGRICULTURE"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 1336.8
$ws1.Range("E8").Value = 338.17
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "BRVM - INDUSTRIE                 (**)"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 3
$ws1.Range("D9").Value = 794.28
$ws1.Range("E9").Value = 266.85
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 685.69
$ws1.Range("E10").Value = 170.24
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM - CONSOMMATION DE BASE          (**)"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 3
$ws1.Range("D11").Value = 664.83
$ws1.Range("E11").Value = 222.65
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 3
$ws1.Range("D12").Value = 656.61
$ws1.Range("E12").Value = 219.41
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

$ws1.Range("A13").Value = "BRVM - FINANCES"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 583.15
$ws1.Range("E13").Value = 145.7
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

$ws1.Range("A14").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 573.1
$ws1.Range("E14").Value = 143.19
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "BRVM-PRESTIGE"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 565.22
$ws1.Range("E15").Value = 141.07
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 493.7
$ws1.Range("E16").Value = 123.77
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

$ws1.Range("A17").Value = "BRVM - ENERGIE"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 447.76
$ws1.Range("E17").Value = 114.13
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 4
$ws1.Range("D18").Value = 378.3
$ws1.Range("E18").Value = 93.91
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

$ws1.Range("A19").Value = "BRVM - INDUSTRIE                  (**)"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 263.26
$ws1.Range("E19").Value = 263.26
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

$ws1.Range("A20").Value = "BRVM-PRINCIPAL                   (**)"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 218.48
$ws1.Range("E20").Value = 218.48
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

$ws1.Range("A21").Value = "BRVM - CONSOMMATION DE BASE         (**)"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 1
$ws1.Range("D21").Value = 217.81
$ws1.Range("E21").Value = 217.81
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

$ws1.Range("A22").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B22").Value = 2
$ws1.Range("C22").Value = 0
$ws1.Range("D22").Value = 14.99
$ws1.Range("E22").Value = 7.5
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"

$ws1.Range("A23").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B23").Value = 2
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 14.72
$ws1.Range("E23").Value = 7.47
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

$ws1.Range("A24").Value = "SAFCA CI (SAFC)"
$ws1.Range("B24").Value = 2
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 9.13
$ws1.Range("E24").Value = -5.58
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "👀 À surveiller"

$ws1.Range("A25").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B25").Value = 3
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = 9.04
$ws1.Range("E25").Value = 7.14
$ws1.Range("F25").Value = "🟢 Achat"
$ws1.Range("G25").Value = "✅ Renforcer"

$ws1.Range("A26").Value = "SICABLE CI (CABC)"
$ws1.Range("B26").Value = 2
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 7.65
$ws1.Range("E26").Value = 7.14
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "👀 À surveiller"

$ws1.Range("A27").Value = "SAPH CI (SPHC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 4.7
$ws1.Range("E27").Value = 4.7
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

$ws1.Range("A28").Value = "SMB CI (SMBC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 4.07
$ws1.Range("E28").Value = 4.07
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = 3.01
$ws1.Range("E29").Value = -2.7
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "👀 À surveiller"

$ws1.Range("A30").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = 2.9
$ws1.Range("E30").Value = -3.08
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "👀 À surveiller"

$ws1.Range("A31").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = 2.24
$ws1.Range("E31").Value = -5.25
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "👀 À surveiller"

$ws1.Range("A32").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = 1.77
$ws1.Range("E32").Value = 4.34
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "👀 À surveiller"

$ws1.Range("A33").Value = "TOTAL"
$ws1.Range("B33").Value = 0
$ws1.Range("C33").Value = 3
$ws1.Range("D33").Value = 0
$ws1.Range("E33").Value = 0
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "➖ Neutre"

$ws1.Range("A34").Value = "BICI CI (BICC)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -0.33
$ws1.Range("E34").Value = 3.85
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "👀 À surveiller"

$ws1.Range("A35").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -1.67
$ws1.Range("E35").Value = -1.67
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "➖ Neutre"

$ws1.Range("A36").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B36").Value = 0
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -2.31
$ws1.Range("E36").Value = -2.31
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "➖ Neutre"

$ws1.Range("A37").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -3.47
$ws1.Range("E37").Value = -3.47
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"

$ws1.Range("A38").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("B38").Value = 0
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -3.85
$ws1.Range("E38").Value = -3.85
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "➖ Neutre"

$ws1.Range("A39").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -3.98
$ws1.Range("E39").Value = -3.98
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "➖ Neutre"

$ws1.Range("A40").Value = "NESTLE CI (NTLC)"
$ws1.Range("B40").Value = 1
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -4.29
$ws1.Range("E40").Value = 3.04
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "👀 À surveiller"

$ws1.Range("A41").Value = "SETAO CI (STAC)"
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -4.35
$ws1.Range("E41").Value = -4.35
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "➖ Neutre"

$ws1.Range("A42").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -4.38
$ws1.Range("E42").Value = -4.38
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"

$ws1.Range("A43").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 2
$ws1.Range("D43").Value = -5.15
$ws1.Range("E43").Value = -3.45
$ws1.Range("F43").Value = "🟡 Observer"
$ws1.Range("G43").Value = "➖ Neutre"

# Table shrank from 44 to 43 data rows (A1:G44 -> A1:G43): drop the last row.
$ws1.Rows.Item(44).Delete() | Out-Null

# --- Sheet "Top_YTD" ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Range("A1").Value = "Titre"
$ws2.Range("B1").Value = "Progression YTD (%)"

$ws2.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B2").Value = 10242702.46

$ws2.Range("A3").Value = "NEI-CEDA CI"
$ws2.Range("B3").Value = 1209044.68

$ws2.Range("A4").Value = "AIR LIQUIDE CI"
$ws2.Range("B4").Value = 391819.1

$ws2.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws2.Range("B5").Value = 228103.76

$ws2.Range("A6").Value = "BRVM - DISTRIBUTION"
$ws2.Range("B6").Value = 130577.17

$ws2.Range("A7").Value = "BRVM - TRANSPORT"
$ws2.Range("B7").Value = 46535.46

$ws2.Range("A8").Value = "BRVM - AGRICULTURE"
$ws2.Range("B8").Value = 35441.19

$ws2.Range("A9").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B9").Value = 5327.17

$ws2.Range("A10").Value = "BRVM - INDUSTRIE                 (**)"
$ws2.Range("B10").Value = 4752.7

$ws2.Range("A11").Value = "BRVM - FINANCES"
$ws2.Range("B11").Value = 3549.55
